# Update the test data strings on the "reg" worksheet (column C, rows 2-11)
# Old -> New:
#   z111     -> z1111
#   Z222     -> Z2222
#   z333     -> z3333
#   z444     -> z4444
#   z555     -> z5555
#   z666     -> z6666
#   z777     -> z7777
#   z888     -> z8888
#   z999     -> z9999
#   z101010  -> z10101010

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

$ws.Range("C2").Value = "z1111"
$ws.Range("C3").Value = "Z2222"
$ws.Range("C4").Value = "z3333"
$ws.Range("C5").Value = "z4444"
$ws.Range("C6").Value = "z5555"
$ws.Range("C7").Value = "z6666"
$ws.Range("C8").Value = "z7777"
$ws.Range("C9").Value = "z8888"
$ws.Range("C10").Value = "z9999"
$ws.Range("C11").Value = "z10101010"
